$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.74203821656051
$ws.Range("C2").Value = 0.691460055096419
$ws.Range("D2").Value = 0.760299625468165
$ws.Range("E2").Value = 0.700934579439252
$ws.Range("F2").Value = 0.616724738675958

$ws.Range("B3").Value = 0.455414012738854
$ws.Range("C3").Value = 0.495867768595041
$ws.Range("D3").Value = 0.49438202247191
$ws.Range("E3").Value = 0.44392523364486
$ws.Range("F3").Value = 0.338850174216028

$ws.Range("B4").Value = 0.522292993630573
$ws.Range("C4").Value = 0.534435261707989
$ws.Range("D4").Value = 0.49063670411985
$ws.Range("E4").Value = 0.47196261682243
$ws.Range("F4").Value = 0.430313588850174

$ws.Range("B5").Value = 0.331210191082803
$ws.Range("C5").Value = 0.421487603305785
$ws.Range("D5").Value = 0.322097378277154
$ws.Range("E5").Value = 0.425233644859813
$ws.Range("F5").Value = 0.376306620209059

$ws.Range("B6").Value = 0.770700636942675
$ws.Range("C6").Value = 0.774104683195592
$ws.Range("D6").Value = 0.741573033707865
$ws.Range("E6").Value = 0.831775700934579
$ws.Range("F6").Value = 0.698606271777004
